# Change the "Project Title" labeled paragraphs so that they just show the
# title text (no "Project Title:" label) in a larger, bold font.
$d = $word.ActiveDocument

$label = "Project Title"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $paraText = $p.Range.Text

    if ($paraText.StartsWith($label)) {
        # Pull out just the title value, e.g. "Project Title: M365 Monthly Release"
        # -> "M365 Monthly Release"
        $rest = $paraText.Substring($label.Length)
        $rest = $rest.TrimEnd([char]13, [char]7)
        if ($rest.StartsWith(":")) {
            $rest = $rest.Substring(1)
        }
        $titleValue = $rest.TrimStart(" ")

        $r = $p.Range
        # Exclude the trailing paragraph mark from the replacement range.
        $newRange = $d.Range($r.Start, $r.End - 1)
        $newRange.Text = $titleValue
        $newRange.Font.Bold = 1
        $newRange.Font.Size = 14
    }
}
